$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row for "KITCHENSCO S.A." right before the existing row 5
# (LINCANGO LUGMANIA SANDY LIZETH), pushing everything below it down by one.
$ws1.Rows.Item(5).Insert()
$ws1.Cells.Item(5, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(5, 2).Value = "KITCHENSCO S.A."
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(5, $c).Value = 0
}

# The totals row (previously row 8) is now row 9. Insert a new data row for
# "VIZUETE GALARZA EDWIN RODRIGO" right before it.
$ws1.Rows.Item(9).Insert()
$ws1.Cells.Item(9, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(9, 2).Value = "VIZUETE GALARZA EDWIN RODRIGO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(9, $c).Value = 0
}

# The totals row is now row 10; the "X de 6" labels need to become "X de 8"
# since there are now 8 data rows instead of 6.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(10, $c)
    $oldText = $cell.Text
    $cell.Value = $oldText.Replace(" de 6", " de 8")
}

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# The totals row (row 9) gets a new data row inserted right before it for
# "VIZUETE GALARZA EDWIN RODRIGO".
$ws2.Rows.Item(9).Insert()
$ws2.Cells.Item(9, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(9, 2).Value = "VIZUETE GALARZA EDWIN RODRIGO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(9, $c).Value = 0
}
